$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 88:89, pushing the existing rows 88:203 down
# to 90:205 (carrying all of their data/styles with them, which reproduces
# the "every later weekly observation shifts down two rows" pattern seen in
# the diff, finishing with the old rows 202:203 reappearing as the brand
# new rows 204:205).
$ws.Rows("88:89").Insert()

# Populate the two freshly inserted rows with the new weekly observation
# (same Mercado / Región / Categoría / Variedad / Origen / Clasificación
# metadata as every other row for this market+product, just a new date and
# new volume/price figures for the "Primera" and "Segunda" quality grades).
$ws.Range("A88").Value = 1
$ws.Range("B88").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C88").Value = "Arica y Parinacota"
$ws.Range("D88").Value = 44495
$ws.Range("E88").Value = 15
$ws.Range("F88").Value = 100112043
$ws.Range("G88").Value = "Pepino ensalada"
$ws.Range("H88").Value = "Sin especificar"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 130
$ws.Range("K88").Value = 5000
$ws.Range("L88").Value = 6000
$ws.Range("M88").Value = 5500
$ws.Range("N88").Value = '$/caja 70 unidades'
$ws.Range("O88").Value = "Región de Arica y Parinacota"
$ws.Range("P88").Value = 79
$ws.Range("Q88").Value = 70
$ws.Range("R88").Value = "Hortaliza"

$ws.Range("A89").Value = 1
$ws.Range("B89").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C89").Value = "Arica y Parinacota"
$ws.Range("D89").Value = 44495
$ws.Range("E89").Value = 15
$ws.Range("F89").Value = 100112043
$ws.Range("G89").Value = "Pepino ensalada"
$ws.Range("H89").Value = "Sin especificar"
$ws.Range("I89").Value = "Segunda"
$ws.Range("J89").Value = 150
$ws.Range("K89").Value = 4000
$ws.Range("L89").Value = 5000
$ws.Range("M89").Value = 4500
$ws.Range("N89").Value = '$/caja 100 unidades'
$ws.Range("O89").Value = "Región de Arica y Parinacota"
$ws.Range("P89").Value = 45
$ws.Range("Q89").Value = 100
$ws.Range("R89").Value = "Hortaliza"
